# wetb/prepost/tests/data/demo_dlc/.../dlc01_demos.xlsx
# "add test coverage for prepost.Simulations.ManTurb64"
#
# Turns row 4 (wsp=10 case) from a "no turbulence" steady-wind case into a
# mann-turbulence (ManTurb64, 100x100 box, 512 grid points) case:
#   - Case id. (D4)   : "..._noturb"  -> "..._s100"
#   - wdir (G4)       : 0   -> 100   (re-purposed as turb box size, per diff)
#   - tu_seed (H4)    : 0   -> 1
#   - Turb base name (J4) : "none" -> "turb_s100_10ms"
#   - turb_dx (K4)    : E4*B4/8192 -> E4*B4/512

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Formula = '="dlc01_steady_wsp" & E4 & "_s100"'
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 1
$ws.Range("J4").Value = "turb_s100_10ms"
$ws.Range("K4").Formula = "=E4*B4/512"

# Restore the active-cell selection recorded in the saved view state.
[void]$ws.Range("Q16").Select()

# Best-effort: the saved window's tab-ratio (horizontal scrollbar/sheet-tab
# splitter position) moved slightly (992 -> 991) in the source file.
$excel.ActiveWindow.TabRatio = 991
